# Updated symbol list on Sat Jan 21 23:37:40 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.52%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.74%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.089"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.20%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07713"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.22%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.192"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-8.98%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.034"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.15%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.019"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.83%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9270"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.94%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09263"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.36%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1823"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.38%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08584"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.76%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "9.38%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09954"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.54%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001477"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.16%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005786"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.32%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.477"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.33%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.14%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3464"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.94%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1325"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.13%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.573"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.58%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2245"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.81%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04679"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.76%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001238"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.44%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004491"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.75%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001307"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.67%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-20.05%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01730"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.61%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04690"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.45%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007887"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.38%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1401"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.73%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007685"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-21.49%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002226"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.23%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008962"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.99%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006220"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.32%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.66%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.787"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "117.98%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002705"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "35.90%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002111"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.66%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002011"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.66%"

